# Eshamy Wild Sockeye - final commit before adage review
# Updates the EVOS pulse/press/pulseRecovery lag columns (G:M) on the DATA
# sheet for rows 21-48 to reflect the recalculated lag values, and restores
# the last active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("DATA")

$cellUpdates = @{
    "K21" = 1
    "L21" = 1
    "M21" = -1
    "H22" = 1
    "I22" = 1
    "J22" = -1
    "L22" = 1
    "M22" = -0.94736842
    "I23" = 1
    "J23" = -0.94736842
    "L23" = 1
    "M23" = -0.89473684
    "H24" = 0
    "J24" = -0.89473684
    "L24" = 1
    "M24" = -0.84210526
    "J25" = -0.84210526
    "K25" = 0
    "M25" = -0.78947368
    "J26" = -0.78947368
    "M26" = -0.73684211
    "J27" = -0.73684211
    "M27" = -0.68421053
    "J28" = -0.68421053
    "M28" = -0.63157895
    "J29" = -0.63157895
    "M29" = -0.57894737
    "J30" = -0.57894737
    "M30" = -0.52631579
    "J31" = -0.52631579
    "M31" = -0.47368421
    "J32" = -0.47368421
    "M32" = -0.42105263
    "J33" = -0.42105263
    "M33" = -0.36842105
    "J34" = -0.36842105
    "M34" = -0.31578947
    "J35" = -0.31578947
    "M35" = -0.26315789
    "J36" = -0.26315789
    "M36" = -0.21052632
    "J37" = -0.21052632
    "M37" = -0.15789474
    "J38" = -0.15789474
    "M38" = -0.10526316
    "J39" = -0.10526316
    "M39" = -0.05263158
    "J40" = -0.05263158
    "M40" = 0
    "J41" = 0
    "M41" = 0
    "J42" = 0
    "M42" = 0
    "G43" = 0
    "M43" = 0
    "G44" = 0
    "J44" = 0
    "G45" = 0
    "J45" = 0
    "M45" = 0
    "G46" = 0
    "J46" = 0
    "M46" = 0
    "G47" = 0
    "J47" = 0
    "M47" = 0
    "G48" = 0
    "J48" = 0
    "M48" = 0
}

foreach ($addr in $cellUpdates.Keys) {
    $ws.Range($addr).Value = $cellUpdates[$addr]
}

$ws.Activate()
[void]$ws.Range("L44").Select()
